# Applies the "Updated cryptos list" price/volume refresh to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.866.25"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "1.814.54"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("D5").Value = "'308.88"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4659"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("D8").Value = "'0.3683"
$ws.Range("E8").Value = "  -0.86%  "

$ws.Range("D9").Value = "'0.07361"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("D11").Value = "'20.39"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "1.837.32"
$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").Value = "'5.366"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.506"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.07059"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").Value = "'91.52"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'0.000008672"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'14.71"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").Value = "26.890.54"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "'5.337"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").Value = "2.039.53"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").Value = "'1.901"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("D26").Value = "'150.27"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").Value = "'2.174"

$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("D29").Value = "'5.321"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").Value = "'115.56"
$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("D31").Value = "'0.08905"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "'0.7665"
$ws.Range("E32").Value = "  +1.76%  "

$ws.Range("D33").Value = "'1.164"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("D34").Value = "'4.508"
$ws.Range("E34").Value = "  +1.77%  "

$ws.Range("D35").Value = "'2.902"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'1.089"
$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("D38").Value = "'0.01959"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").Value = "'0.05279"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").Value = "'7.246"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").Value = "'0.5317"
$ws.Range("E42").Value = "  +1.83%  "

$ws.Range("D43").Value = "'2.344"
$ws.Range("E43").Value = "  -1.11%  "

$ws.Range("D44").Value = "'0.1659"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "'8.424"
$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("D46").Value = "'0.4924"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "'1.670"
$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("D50").Value = "'103.75"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").Value = "'0.06284"
$ws.Range("E51").Value = "  +0.00%  "
